$d = $word.ActiveDocument

# =======================================================================
# 1) "The system must notify the user when it was last active"
#    -> "The system must notify the user that it is active."
#    (only the first occurrence of this phrase in the document; the
#    later "... and when it was last active." inside the
#    notifyStartUp.py description is left untouched)
#
#    A short-lived bookmark is dropped right at the run boundary before
#    the replacement; that keeps the untouched " the user " run from
#    being silently coalesced into the run being edited, then the
#    scratch bookmark is removed again once the text is in place.
# =======================================================================
$rngFind1 = $d.Content
$found1 = $rngFind1.Find.Execute("when it was last active", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $boundary1 = $d.Range($rngFind1.Start, $rngFind1.Start)
    $d.Bookmarks.Add("zzScratch1", $boundary1) | Out-Null

    $rngReplace1 = $d.Content
    $rngReplace1.Find.Execute("when it was last active", $true, $false, $false, $false, $false, $true, 1, $false, "that it is active.", 1) | Out-Null

    $d.Bookmarks("zzScratch1").Delete()
} else {
    Write-Output "WARNING: first replacement target not found"
}

# =======================================================================
# 2) "... via an e-mail containing the last streamed video and the
#    timestamp of when it was last active."
#    -> "... via an e-mail containing the last streamed video and the
#    timestamp." with the "_GoBack" bookmark now sitting right before
#    the new final period.
#
#    Same scratch-bookmark trick keeps "The user is notified" from
#    merging into the run being shortened. Once the text reads "...the
#    timestamp." the real "_GoBack" bookmark is (re)created right
#    before the trailing period. A document only ever has a single
#    "_GoBack" bookmark, so this also removes it from its old spot at
#    the very start of the document (the "Surveillance System"
#    heading).
# =======================================================================
$rngFind2 = $d.Content
$found2 = $rngFind2.Find.Execute(" via an e-mail containing the last streamed video and the timestamp of when it was last active.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $boundary2 = $d.Range($rngFind2.Start, $rngFind2.Start)
    $d.Bookmarks.Add("zzScratch2", $boundary2) | Out-Null

    $rngReplace2 = $d.Content
    $rngReplace2.Find.Execute(" via an e-mail containing the last streamed video and the timestamp of when it was last active.", $true, $false, $false, $false, $false, $true, 1, $false, " via an e-mail containing the last streamed video and the timestamp.", 1) | Out-Null

    $d.Bookmarks("zzScratch2").Delete()

    $rngPeriod = $d.Content
    $foundPeriod = $rngPeriod.Find.Execute("the timestamp.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundPeriod) {
        $periodPos = $rngPeriod.End - 1
        $bookmarkPoint = $d.Range($periodPos, $periodPos)
        $d.Bookmarks.Add("_GoBack", $bookmarkPoint) | Out-Null
    } else {
        Write-Output "WARNING: could not locate insertion point for _GoBack bookmark"
    }
} else {
    Write-Output "WARNING: second replacement target not found"
}
